# Update "想去人数" (F column) figures in both the "展览" and "全部类型"
# sheets to reflect refreshed output data.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row 4 (F4: 810 -> 813), row 6 (F6: 23 -> 24)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 813
$wsExhibit.Range("F6").Value = 24

# Sheet "全部类型": row 5 (F5: 810 -> 813), row 7 (F7: 23 -> 24)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 813
$wsAll.Range("F7").Value = 24
